# Daily attendance processing - reorder the "Recorded By" (column G) author
# lists for specific known value-combinations by reversing their order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Known "Recorded By" combinations that need to be reordered (reversed),
# mapped from their current text to the new text.
$map = @{
    "System, dnasr281@gmail.com" = "dnasr281@gmail.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
    "backup@backdoor.com, System" = "System, backup@backdoor.com"
    "backup@backdoor.com, System, system" = "system, System, backup@backdoor.com"
}

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
